# Add new power plant types to the "Electricity Source" subscript list
# on the FSCaFoCC sheet (issues #280 and #99), and tidy up the
# formatting of the existing subscript labels (remove stray fill from
# the label column, and clear the redundant format on B1).

$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("FSCaFoCC")

# --- Normalize formatting of the existing rows -----------------------
# A2:A18 used a style that redundantly applied fill; make them match
# the already-present bold/no-fill style used elsewhere on the sheet.
$ws.Range("A2:A18").Font.Bold = $true

# B1 used a style with redundant (no-op) fill/alignment flags; clear it
# back to the default style.
$ws.Range("B1").ClearFormats()

# --- Append the new Electricity Source entries ------------------------
$newPlantTypes = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$row = 19
foreach ($plantType in $newPlantTypes) {
    $ws.Cells.Item($row, 1).Value2 = $plantType
    $ws.Cells.Item($row, 2).Value2 = 0
    $row++
}

# Match the bold styling used by the rest of the subscript labels.
$ws.Range("A19:A24").Font.Bold = $true

# --- Restore the view state -------------------------------------------
# Reflect that the FSCaFoCC sheet was last scrolled/selected at the new
# blank row (A25) before switching focus back to the About sheet.
$ws.Activate() | Out-Null
$ws.Range("A25").Select() | Out-Null
$wsAbout.Activate() | Out-Null
